$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '21.663.85'
$ws.Cells.Item(2, 5).Value = '  -1.57%  '
$ws.Cells.Item(3, 4).Value = '1.533.02'
$ws.Cells.Item(3, 5).Value = '  -1.42%  '
$ws.Cells.Item(4, 5).Value = '  +0.10%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '288.29'
$ws.Cells.Item(6, 5).Value = '  +0.46%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.3942'
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3153'
$ws.Cells.Item(8, 5).Value = '  -2.49%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '42.09'
$ws.Cells.Item(9, 5).Value = '  +1.78%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.07154'
$ws.Cells.Item(10, 5).Value = '  -2.14%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '1.043'
$ws.Cells.Item(11, 5).Value = '  -6.78%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '1.002'
$ws.Cells.Item(12, 5).Value = '  +0.21%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '5.627'
$ws.Cells.Item(13, 5).Value = '  -1.55%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '18.45'
$ws.Cells.Item(14, 5).Value = '  -4.70%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '6.602'
$ws.Cells.Item(15, 5).Value = '  -3.08%  '
$ws.Cells.Item(16, 4).Value = '1.536.62'
$ws.Cells.Item(16, 5).Value = '  -1.24%  '
$ws.Cells.Item(17, 5).Value = '  -0.13%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.06591'
$ws.Cells.Item(18, 5).Value = '  -0.56%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '82.90'
$ws.Cells.Item(19, 5).Value = '  -2.79%  '
$ws.Cells.Item(20, 5).Value = '  +0.14%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '6.088'
$ws.Cells.Item(21, 5).Value = '  -4.96%  '
$ws.Cells.Item(22, 5).Value = '  -3.48%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '10.81'
$ws.Cells.Item(23, 5).Value = '  -5.48%  '
$ws.Cells.Item(24, 5).Value = '  +4.22%  '
$ws.Cells.Item(25, 4).Value = '21.659.51'
$ws.Cells.Item(25, 5).Value = '  -1.59%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '2.328'
$ws.Cells.Item(26, 5).Value = '  -7.66%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '147.43'
$ws.Cells.Item(27, 5).Value = '  -0.92%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '18.28'
$ws.Cells.Item(28, 5).Value = '  -2.83%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '4.839'
$ws.Cells.Item(29, 5).Value = '  -0.38%  '
$ws.Cells.Item(30, 4).Value = '1.711.53'
$ws.Cells.Item(30, 5).Value = '  -0.99%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '116.85'
$ws.Cells.Item(31, 5).Value = '  -3.04%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '5.839'
$ws.Cells.Item(32, 5).Value = '  -0.55%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.9467'
$ws.Cells.Item(33, 5).Value = '  -13.30%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.08125'
$ws.Cells.Item(34, 5).Value = '  -0.12%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '8.497'
$ws.Cells.Item(35, 5).Value = '  -8.32%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.06049'
$ws.Cells.Item(36, 5).Value = '  -2.53%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '5.076'
$ws.Cells.Item(37, 5).Value = '  -3.19%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.02199'
$ws.Cells.Item(38, 5).Value = '  -4.25%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '1.439'
$ws.Cells.Item(39, 5).Value = '  -12.33%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.2013'
$ws.Cells.Item(40, 5).Value = '  -4.31%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.176'
$ws.Cells.Item(41, 5).Value = '  -3.43%  '
$ws.Cells.Item(42, 5).Value = '  +0.10%  '
$ws.Cells.Item(43, 5).Value = '  -0.34%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.5728'
$ws.Cells.Item(44, 5).Value = '  -3.32%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '12.97'
$ws.Cells.Item(45, 5).Value = '  -3.77%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '3.720'
$ws.Cells.Item(46, 5).Value = '  -0.04%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.5466'
$ws.Cells.Item(47, 5).Value = '  -4.71%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '1.151'
$ws.Cells.Item(48, 5).Value = '  -0.10%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.857'
$ws.Cells.Item(49, 5).Value = '  -3.79%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '115.26'
$ws.Cells.Item(50, 5).Value = '  -3.40%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.06686'
$ws.Cells.Item(51, 5).Value = '  -2.64%  '
